$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "關鍵字" (keywords) column C, mirroring the 問題/答案 (question/answer) data
$ws.Range('C1').Value = '關鍵字'

$keywords = @(
    '109學年 吳文傑老師 全英授課 經濟學課程',
    '這學期 早八 社會學',
    '111學年第二學期 吳岳剛老師 靜態影像設計 基礎攝影技巧',
    '幽默 實用 生活化 民法概要',
    '穆斯林 種族議題 課',
    '企業模式 實務課程',
    '輕鬆學 印尼文 課',
    '金融市場 銀行 運作模式 一學年 貨幣銀行學',
    '晚九 口頭報告 營養學分',
    '台灣大車隊 副董事長 管理學 課程',
    '政大 課程 數位學術期刊',
    '民意調查 課',
    '政大 民族系 民族文學 課',
    '資管系 系學會 服務學習 課程',
    '人口遷移測量 課',
    '政大 林怡伶 程式設計一 課程 難度',
    '政大 林巧敏 圖書館管理 課程 內容 難度',
    '政大 吳漢銘 統計學（一） 課程 難度',
    '政大 張仲豪 體育課 武術初級 基礎 推薦',
    '111-1 政大 吳豐祥 科技與創新管理 課 內容 難度',
    '110-2 政大 林宛瑩 財務報表分析 課 老師 嚴格 點名',
    '111-2 政大 朱斌妤 溝通與協商 課程 互動 知識 分數',
    '111-1 政大 甯方璽 測量學及實習 課 內容 考試 分數',
    '110-2 政大 陳雅莉 大學英文（二） 課 答問題 報告 長篇文章',
    '111-1 政大 林信助 經濟學 課 點名 功課 問題 拍影片 報告',
    '110-2 政大 黃柏鈞 個體經濟學 課 一般經濟學 差異 難度',
    '109-1 政大 趙知章 生活中的生命科學 通識課 助教 點名 考試 普通生物 知識',
    '109-2 政大 詹凌菁 高級會計學（二） 課 一般會計學 功課 考試 點名',
    '111-2 政大 陳立夫 土地行政 課程 內容 深度 實地勘察 報告',
    '110-2 政大 林文乙 橄欖球初級體育課 曬太陽 容易受傷 瘦弱 适合'
)

for ($i = 0; $i -lt $keywords.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $keywords[$i]
}

# Match the formatting of the existing columns instead of leaving the default style
# (copy the header style from A1 onto the new header cell, and the body style
# from B2 onto the whole new data range, same as the rest of column B).
$ws.Range('A1').Copy() | Out-Null
$ws.Range('C1').PasteSpecial(-4122) | Out-Null
$ws.Range('B2').Copy() | Out-Null
$ws.Range('C2:C31').PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Selection ends up on D5 (no more frozen top-left scroll position on A17/A30)
$ws.Range('D5').Select() | Out-Null

